# Re-process the metadata: columns G (nivel-formativo-grupo-iaest-descripcion)
# and H (sexo) move from being "dimension" columns to "measure" columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-dimension:X -> iaest-measure:X
$ws.Range("G2").Value = "iaest-measure:nivel-formativo-grupo-iaest-descripcion"
$ws.Range("H2").Value = "iaest-measure:sexo"

# Row 3: dim -> medida
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"

# Row 4: skos:Concept -> xsd:int
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"

# Row 5: the mapping-file cells no longer apply to measures, remove them entirely
$ws.Range("G5").Clear() | Out-Null
$ws.Range("H5").Clear() | Out-Null

Write-Host "edit applied"
